$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-25 22:18:37"
$ws.Range("H2").Value = "'48%"
$ws.Range("I2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = "2026-02-25 22:18:40"
$ws.Range("N3").Value = "0.8 °C 21:51 TU"
$ws.Range("E4").Value = "2026-02-25 22:18:42"
$ws.Range("J4").Value = "1022.1 hPa"
$ws.Range("E5").Value = "2026-02-25 22:18:45"
$ws.Range("N5").Value = "1.5 °C 21:32 TU"
$ws.Range("O5").Value = "5.5 °C"
$ws.Range("E6").Value = "2026-02-25 22:18:48"
$ws.Range("J6").Value = "1022.0 hPa"
$ws.Range("E7").Value = "2026-02-25 22:18:50"
$ws.Range("J7").Value = "1021.6 hPa"
$ws.Range("E8").Value = "2026-02-25 22:18:53"
$ws.Range("J8").Value = "1021.3 hPa"
$ws.Range("E9").Value = "2026-02-25 22:18:56"
$ws.Range("E10").Value = "2026-02-25 22:18:58"
$ws.Range("O10").Value = "9.2 °C"
$ws.Range("E11").Value = "2026-02-25 22:19:01"
$ws.Range("H11").Value = "'64%"
$ws.Range("I11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("O11").Value = "8.7 °C"
$ws.Range("E12").Value = "2026-02-25 22:19:03"
$ws.Range("E13").Value = "2026-02-25 22:19:06"
$ws.Range("H13").Value = "'65%"
$ws.Range("I13").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("J13").Value = "1023.0 hPa"
$ws.Range("E14").Value = "2026-02-25 22:19:08"
$ws.Range("E15").Value = "2026-02-25 22:19:11"
$ws.Range("E16").Value = "2026-02-25 22:19:14"
$ws.Range("H16").Value = "'33%"
$ws.Range("I16").Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = "2026-02-25 22:19:16"
$ws.Range("O17").Value = "9.0 °C"
$ws.Range("E18").Value = "2026-02-25 22:19:19"
$ws.Range("H18").Value = "'90%"
$ws.Range("I18").Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4122) | Out-Null
$ws.Range("J18").Value = "1022.2 hPa"
$ws.Range("E19").Value = "2026-02-25 22:19:21"
$ws.Range("N19").Value = "7.5 °C 21:59 TU"
$ws.Range("O19").Value = "12.0 °C"
$ws.Range("E20").Value = "2026-02-25 22:19:24"
$ws.Range("N20").Value = "-1.9 °C 21:39 TU"
$ws.Range("O20").Value = "2.6 °C"
$ws.Range("E21").Value = "2026-02-25 22:19:27"
$ws.Range("J21").Value = "1021.7 hPa"
$ws.Range("O21").Value = "10.1 °C"
$ws.Range("E22").Value = "2026-02-25 22:19:29"
$ws.Range("H22").Value = "'45%"
$ws.Range("I22").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "2026-02-25 22:19:32"
$ws.Range("H23").Value = "'36%"
$ws.Range("I23").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4122) | Out-Null
$ws.Range("K23").Value = "16.4 MJ/m2"
$ws.Range("O23").Value = "3.8 °C"
$ws.Range("E24").Value = "2026-02-25 22:19:34"
$ws.Range("J24").Value = "1020.3 hPa"
$ws.Range("L24").Value = "21.2 km/h - 88º 21:31 TU"
$ws.Range("E25").Value = "2026-02-25 22:19:37"
$ws.Range("H25").Value = "'35%"
$ws.Range("I25").Copy() | Out-Null
$ws.Range("H25").PasteSpecial(-4122) | Out-Null
$ws.Range("N25").Value = "1.5 °C 21:59 TU"
$ws.Range("O25").Value = "5.1 °C"
$ws.Range("E26").Value = "2026-02-25 22:19:40"
$ws.Range("E27").Value = "2026-02-25 22:19:42"
$ws.Range("N27").Value = "1.3 °C 21:47 TU"
$ws.Range("O27").Value = "5.1 °C"
$ws.Range("E28").Value = "2026-02-25 22:19:45"
$ws.Range("H28").Value = "'85%"
$ws.Range("I28").Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$ws.Range("J28").Value = "1022.0 hPa"
$ws.Range("O28").Value = "8.9 °C"
$ws.Range("E29").Value = "2026-02-25 22:19:48"
$ws.Range("E30").Value = "2026-02-25 22:19:50"
$ws.Range("H30").Value = "'92%"
$ws.Range("I30").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null
$ws.Range("E31").Value = "2026-02-25 22:19:53"
$ws.Range("E32").Value = "2026-02-25 22:19:55"
$ws.Range("O32").Value = "9.0 °C"
$ws.Range("E33").Value = "2026-02-25 22:19:58"
$ws.Range("J33").Value = "1021.5 hPa"
$ws.Range("O33").Value = "8.4 °C"
$ws.Range("E34").Value = "2026-02-25 22:20:01"
$ws.Range("H34").Value = "'54%"
$ws.Range("I34").Copy() | Out-Null
$ws.Range("H34").PasteSpecial(-4122) | Out-Null
$ws.Range("E35").Value = "2026-02-25 22:20:03"
$ws.Range("G35").Value = "2 cm"
$ws.Range("J35").Value = "1019.8 hPa"
$ws.Range("E36").Value = "2026-02-25 22:20:06"
$ws.Range("E37").Value = "2026-02-25 22:20:09"
$ws.Range("H37").Value = "'84%"
$ws.Range("I37").Copy() | Out-Null
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("J37").Value = "1023.7 hPa"
$ws.Range("E38").Value = "2026-02-25 22:20:11"
$ws.Range("O38").Value = "9.3 °C"
$ws.Range("E39").Value = "2026-02-25 22:20:14"
$ws.Range("E40").Value = "2026-02-25 22:20:16"
$ws.Range("J40").Value = "1022.1 hPa"
$ws.Range("O40").Value = "9.4 °C"
$ws.Range("E41").Value = "2026-02-25 22:20:19"
$ws.Range("J41").Value = "1021.2 hPa"
$ws.Range("E42").Value = "2026-02-25 22:20:22"
$ws.Range("E43").Value = "2026-02-25 22:20:24"
$ws.Range("O43").Value = "9.7 °C"
$ws.Range("E44").Value = "2026-02-25 22:20:27"
$ws.Range("H44").Value = "'47%"
$ws.Range("I44").Copy() | Out-Null
$ws.Range("H44").PasteSpecial(-4122) | Out-Null
$ws.Range("O44").Value = "2.1 °C"
$ws.Range("E45").Value = "2026-02-25 22:20:29"
$ws.Range("J45").Value = "1020.1 hPa"
$ws.Range("E46").Value = "2026-02-25 22:20:32"
$ws.Range("H46").Value = "'82%"
$ws.Range("I46").Copy() | Out-Null
$ws.Range("H46").PasteSpecial(-4122) | Out-Null
$ws.Range("J46").Value = "1021.0 hPa"
$ws.Range("O46").Value = "9.9 °C"
